$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - remove all existing content so the shared string table
# is rebuilt from scratch in the order we set values below.
$ws.Cells.Clear()

# Text cells - set in this order so the rebuilt shared-strings table
# lines up with the workbook produced by the "Excel Dump" tool.
$ws.Range("A2").Value = "name"
$ws.Range("D2").Value = "row1"
$ws.Range("D6").Value = "s"
$ws.Range("A1").Value = "tkb_testing_scraptemp"
$ws.Range("B2").Value = "Category"
$ws.Range("C2").Value = "Cost"

# Data rows
$ws.Range("A3").Value = 6686
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 3.45
$ws.Range("D3").Value = 2

$ws.Range("A4").Value = 6729
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 2.75
$ws.Range("D4").Value = 2

$ws.Range("A5").Value = 4900
$ws.Range("B5").Value = 30
$ws.Range("C5").Value = 2.44

$ws.Range("A6").Value = 4916
$ws.Range("B6").Value = 40
$ws.Range("C6").Value = 1.88

$ws.Range("D3").Select()
